{"js": "// Homework1.docx date bump:\n//   Release date: Mar 21 -> Nov 17\n//   Due date:     Apr 11 -> Dec 1\n//\n// Each date lives in its own paragraph (\"Release date: Mar 21, 2020\" /\n// \"Due date: Apr 11, 2020 23:59\"), each split across several runs\n// (month / space / day / \", 20\" / year...). We scope every search to the\n// specific paragraph so short tokens like \"11\" or \"21\" (which also show\n// up elsewhere in the document, e.g. reference numbers \"[11]\") are never\n// ambiguously matched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the \"Release date:\" and \"Due date:\" paragraphs by content so the\n// script doesn't depend on a brittle fixed paragraph index.\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst releaseParagraph = paragraphs.items.find((p) => p.text.indexOf(\"Release date:\") !== -1);\nconst dueParagraph = paragraphs.items.find((p) => p.text.indexOf(\"Due date:\") !== -1);\n\nif (!releaseParagraph) {\n  throw new Error(\"Could not find the 'Release date:' paragraph.\");\n}\nif (!dueParagraph) {\n  throw new Error(\"Could not find the 'Due date:' paragraph.\");\n}\n\n// --- Release date: \"Mar\" -> \"Nov\", \"21\" -> \"17\" -------------------------\nconst monthResults1 = releaseParagraph.search(\"Mar\", { matchCase: true, matchWholeWord: true });\nmonthResults1.load(\"items\");\nawait context.sync();\nmonthResults1.items[0].insertText(\"Nov\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst dayResults1 = releaseParagraph.search(\"21\", { matchCase: true, matchWholeWord: true });\ndayResults1.load(\"items\");\nawait context.sync();\ndayResults1.items[0].insertText(\"17\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Due date: \"Apr\" -> \"Dec\", \"11\" -> \"1\" -------------------------------\nconst monthResults2 = dueParagraph.search(\"Apr\", { matchCase: true, matchWholeWord: true });\nmonthResults2.load(\"items\");\nawait context.sync();\nmonthResults2.items[0].insertText(\"Dec\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst dayResults2 = dueParagraph.search(\"11\", { matchCase: true, matchWholeWord: true });\ndayResults2.load(\"items\");\nawait context.sync();\ndayResults2.items[0].insertText(\"1\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Homework1.docx date bump:\n#   Release date: Mar 21 -> Nov 17\n#   Due date:     Apr 11 -> Dec 1\n#\n# Short tokens like \"11\"/\"21\" also appear elsewhere in the document (e.g.\n# reference \"[11]\"), so every Find is scoped to the specific \"Release\n# date:\" / \"Due date:\" paragraph (located by its own text) instead of\n# running over the whole document.\n\n$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$releasePara = $null\n$duePara = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"Release date:*\") {\n        $releasePara = $para\n    }\n    if ($para.Range.Text -like \"Due date:*\") {\n        $duePara = $para\n    }\n}\n\nif ($releasePara -eq $null) {\n    throw \"Could not find the 'Release date:' paragraph.\"\n}\nif ($duePara -eq $null) {\n    throw \"Could not find the 'Due date:' paragraph.\"\n}\n\n# --- Release date: \"Mar\" -> \"Nov\", \"21\" -> \"17\" --------------------------\n$rng = $releasePara.Range\n$rng.Find.Execute(\"Mar\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $null, \"Nov\", $wdReplaceOne)\n\n$rng = $releasePara.Range\n$rng.Find.Execute(\"21\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $null, \"17\", $wdReplaceOne)\n\n# --- Due date: \"Apr\" -> \"Dec\", \"11\" -> \"1\" --------------------------------\n$rng = $duePara.Range\n$rng.Find.Execute(\"Apr\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $null, \"Dec\", $wdReplaceOne)\n\n$rng = $duePara.Range\n$rng.Find.Execute(\"11\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $null, \"1\", $wdReplaceOne)\n"}
